$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force the value to be stored as text (shared string), matching the
    # original workbook's convention of keeping all values - even
    # numeric-looking ones - as text cells.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $wsLider.Range("A2") "0.049999999999998934 - x + y"
Set-TextValue $wsLider.Range("B2") "-0.049999999999998934"
Set-TextValue $wsLider.Range("D2") "0.4"

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $wsFollower.Range("A2") "-28.85 + x + y"
Set-TextValue $wsFollower.Range("B2") "8.85"
Set-TextValue $wsFollower.Range("D2") "0.55"
Set-TextValue $wsFollower.Range("F2") "2.8000000000000003"

Set-TextValue $wsFollower.Range("A3") "14.4 - y"
Set-TextValue $wsFollower.Range("B3") "-14.4"
Set-TextValue $wsFollower.Range("D3") "0.45"
Set-TextValue $wsFollower.Range("E3") "-4.6000000000000005"
Set-TextValue $wsFollower.Range("F3") "-3.4000000000000004"

Set-TextValue $wsFollower.Range("A4") "-25.6 + y"
Set-TextValue $wsFollower.Range("B4") "-5.6"
Set-TextValue $wsFollower.Range("D4") "0.55"
Set-TextValue $wsFollower.Range("E4") "2.3000000000000003"
Set-TextValue $wsFollower.Range("F4") "0"

$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "14.45"
Set-TextValue $wsPunto.Range("B2") "14.4"

# NOTE: the sheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively in this environment, so
# name-based lookup would make both variables alias the very same sheet.
# Use the (1-based) sheet index instead, matching the workbook's sheet order:
# 1 Funciones_Objetivo, 2 Restricciones_del_lider, 3 Restricciones_del_follower,
# 4 Punto_modificado, 5 Vector_bf, 6 Vector_BF, 7 Vector_Alpha.
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf.Range("A2") "-53.65"

$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "-28.5"
Set-TextValue $wsBF.Range("A3") "-16.1"
